# 1. Table on slide 16 switches from the custom "Table_0" style to the
#    built-in PowerPoint table style {189BC292-AAF2-4C9C-B82D-26BBCE7D9268}.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(16)
$tbl = $s.Shapes.Item(3).Table
$tbl.ApplyStyle("{189BC292-AAF2-4C9C-B82D-26BBCE7D9268}")

# 2. The presentation's theme colour scheme (the "Integral" design) is
#    swapped for the stock "Office Theme" colours. All slides share the
#    one design/master, so editing it through slide 1 re-writes the
#    underlying theme part used everywhere.
$first = $p.Slides.Item(1)
$colors = $first.ThemeColorScheme
$colors.Item(1).RGB  = 0         # dk1      000000
$colors.Item(2).RGB  = 16777215  # lt1      FFFFFF
$colors.Item(3).RGB  = 6968388   # dk2      44546A
$colors.Item(4).RGB  = 15132391  # lt2      E7E6E6
$colors.Item(5).RGB  = 13998939  # accent1  5B9BD5
$colors.Item(6).RGB  = 3243501   # accent2  ED7D31
$colors.Item(7).RGB  = 10855845  # accent3  A5A5A5
$colors.Item(8).RGB  = 49407     # accent4  FFC000
$colors.Item(9).RGB  = 12874308  # accent5  4472C4
$colors.Item(10).RGB = 4697456   # accent6  70AD47
$colors.Item(11).RGB = 12673797  # hlink    0563C1
$colors.Item(12).RGB = 7491477   # folHlink 954F72
